$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 743.1667
$ws.Range("I17").Value = 1119.2307
$ws.Range("J17").Value = 650.9245
$ws.Range("K17").Value = 3357.6921
$ws.Range("L17").Value = 1952.7735
$ws.Range("M17").Value = -3189.6921
$ws.Range("N17").Value = -2288.7735
$ws.Range("H75").Value = 39315
$ws.Range("J75").Value = 39315
$ws.Range("L75").Value = 39315
$ws.Range("N75").Value = -41187
$ws.Range("H78").Value = 39315
$ws.Range("J78").Value = 39315
$ws.Range("L78").Value = 117945
$ws.Range("N78").Value = -127305
$ws.Range("H128").Value = 41612.5
$ws.Range("J128").Value = 41612.5
$ws.Range("L128").Value = 41612.5
$ws.Range("N128").Value = -51572.5
$ws.Range("H132").Value = 38468880
$ws.Range("I132").Value = 62509290
$ws.Range("J132").Value = 4230.5
$ws.Range("K132").Value = 187527870
$ws.Range("L132").Value = 12691.5
$ws.Range("M132").Value = -187525340
$ws.Range("N132").Value = -17751.5
$ws.Range("H135").Value = 1022
$ws.Range("I135").Value = 630.8
$ws.Range("K135").Value = 5677.2
$ws.Range("M135").Value = -3142.2
$ws.Range("H137").Value = 1895.2106
$ws.Range("I137").Value = 937.6957
$ws.Range("J137").Value = 5899.364
$ws.Range("K137").Value = 2813.0871
$ws.Range("L137").Value = 17698.092
$ws.Range("M137").Value = -263.0870999999997
$ws.Range("N137").Value = -22798.092
$ws.Range("H138").Value = 3986.5
$ws.Range("I138").Value = 1156.9546
$ws.Range("J138").Value = 4805.579
$ws.Range("K138").Value = 3470.8638
$ws.Range("L138").Value = 14416.737
$ws.Range("M138").Value = 1669.1362
$ws.Range("N138").Value = -24696.737
$ws.Range("H141").Value = 12506.947
$ws.Range("I141").Value = 13403.059
$ws.Range("J141").Value = 4890
$ws.Range("K141").Value = 40209.177
$ws.Range("L141").Value = 14670
$ws.Range("M141").Value = -35029.177
$ws.Range("N141").Value = -25030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4454.0684
$ws.Range("I32").Value = 3300.5
$ws.Range("J32").Value = 15330.571
$ws.Range("K32").Value = 3300.5
$ws.Range("L32").Value = 15330.571
$ws.Range("M32").Value = -3013.5
$ws.Range("N32").Value = -15904.571
$ws.Range("H61").Value = 1100.8667
$ws.Range("I61").Value = 949.5263
$ws.Range("J61").Value = 1362.2727
$ws.Range("K61").Value = 949.5263
$ws.Range("L61").Value = 1362.2727
$ws.Range("M61").Value = -737.5263
$ws.Range("N61").Value = -1786.2727
$ws.Range("H74").Value = 3001.975
$ws.Range("I74").Value = 3139.258
$ws.Range("J74").Value = 2529.111
$ws.Range("K74").Value = 3139.258
$ws.Range("L74").Value = 2529.111
$ws.Range("M74").Value = -2265.258
$ws.Range("N74").Value = -4277.111
$ws.Range("H77").Value = 3001.975
$ws.Range("I77").Value = 3139.258
$ws.Range("J77").Value = 2529.111
$ws.Range("K77").Value = 15696.29
$ws.Range("L77").Value = 12645.555
$ws.Range("M77").Value = -11328.29
$ws.Range("N77").Value = -21381.555
$ws.Range("H136").Value = 1100.8667
$ws.Range("I136").Value = 949.5263
$ws.Range("J136").Value = 1362.2727
$ws.Range("K136").Value = 2848.5789
$ws.Range("L136").Value = 4086.8181
$ws.Range("M136").Value = -298.5789
$ws.Range("N136").Value = -9186.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2549.6
$ws.Range("I107").Value = 2474
$ws.Range("J107").Value = 2663
$ws.Range("K107").Value = 2474
$ws.Range("L107").Value = 2663
$ws.Range("M107").Value = -554
$ws.Range("N107").Value = -6503

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1957.1449
$ws.Range("I58").Value = 1683.0159
$ws.Range("K58").Value = 1683.0159
$ws.Range("M58").Value = -1480.0159
$ws.Range("H106").Value = 34833.332
$ws.Range("J106").Value = 34833.332
$ws.Range("L106").Value = 34833.332
$ws.Range("N106").Value = -37357.332
$ws.Range("H115").Value = 29700
$ws.Range("J115").Value = 29700
$ws.Range("L115").Value = 29700
$ws.Range("N115").Value = -32050
$ws.Range("H132").Value = 5033.5
$ws.Range("I132").Value = 4592.615
$ws.Range("J132").Value = 6179.8
$ws.Range("K132").Value = 13777.845
$ws.Range("L132").Value = 18539.4
$ws.Range("M132").Value = -11247.845
$ws.Range("N132").Value = -23599.4
$ws.Range("H134").Value = 8917.764999999999
$ws.Range("I134").Value = 21180.4
$ws.Range("J134").Value = 3808.3333
$ws.Range("K134").Value = 63541.2
$ws.Range("L134").Value = 11424.9999
$ws.Range("M134").Value = -61006.2
$ws.Range("N134").Value = -16494.9999
$ws.Range("H136").Value = 1957.1449
$ws.Range("I136").Value = 1683.0159
$ws.Range("K136").Value = 5049.0477
$ws.Range("M136").Value = -2499.0477

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1495.4828
$ws.Range("I5").Value = 494
$ws.Range("J5").Value = 1946.15
$ws.Range("K5").Value = 1482
$ws.Range("L5").Value = 5838.450000000001
$ws.Range("M5").Value = -1370
$ws.Range("N5").Value = -6062.450000000001
$ws.Range("H62").Value = 7208
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 8587.556
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 25762.668
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -27134.668
$ws.Range("H65").Value = 7208
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 8587.556
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 77288.004
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -84152.004
$ws.Range("H121").Value = 2602.6072
$ws.Range("J121").Value = 2691.037
$ws.Range("L121").Value = 8073.110999999999
$ws.Range("N121").Value = -10693.111
$ws.Range("H131").Value = 6579749.5
$ws.Range("J131").Value = 833.8823
$ws.Range("L131").Value = 2501.6469
$ws.Range("N131").Value = -12581.6469
$ws.Range("H133").Value = 3433.125
$ws.Range("I133").Value = 4632.857
$ws.Range("J133").Value = 2500
$ws.Range("K133").Value = 13898.571
$ws.Range("L133").Value = 7500
$ws.Range("M133").Value = -8838.571
$ws.Range("N133").Value = -17620
$ws.Range("H135").Value = 1495.4828
$ws.Range("I135").Value = 494
$ws.Range("J135").Value = 1946.15
$ws.Range("K135").Value = 4446
$ws.Range("L135").Value = 17515.35
$ws.Range("M135").Value = -1911
$ws.Range("N135").Value = -22585.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 25499.5
$ws.Range("J93").Value = 25499.5
$ws.Range("L93").Value = 25499.5
$ws.Range("N93").Value = -29243.5
$ws.Range("H122").Value = 4759.3
$ws.Range("I122").Value = 2558.6
$ws.Range("J122").Value = 6960
$ws.Range("K122").Value = 7675.799999999999
$ws.Range("L122").Value = 20880
$ws.Range("M122").Value = -5225.799999999999
$ws.Range("N122").Value = -25780
$ws.Range("H132").Value = 2578.8647
$ws.Range("I132").Value = 1569.8
$ws.Range("J132").Value = 4681.0835
$ws.Range("K132").Value = 4709.4
$ws.Range("L132").Value = 14043.2505
$ws.Range("M132").Value = -2179.4
$ws.Range("N132").Value = -19103.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 26900
$ws.Range("J74").Value = 26900
$ws.Range("L74").Value = 26900
$ws.Range("N74").Value = -28896
$ws.Range("H75").Value = 40399.8
$ws.Range("J75").Value = 49999.75
$ws.Range("L75").Value = 49999.75
$ws.Range("N75").Value = -51871.75
$ws.Range("H77").Value = 26900
$ws.Range("J77").Value = 26900
$ws.Range("L77").Value = 80700
$ws.Range("N77").Value = -90684
$ws.Range("H78").Value = 40399.8
$ws.Range("J78").Value = 49999.75
$ws.Range("L78").Value = 149999.25
$ws.Range("N78").Value = -159359.25
$ws.Range("H80").Value = 39750
$ws.Range("J80").Value = 39750
$ws.Range("L80").Value = 39750
$ws.Range("N80").Value = -41996
$ws.Range("H81").Value = 39750
$ws.Range("J81").Value = 39750
$ws.Range("L81").Value = 39750
$ws.Range("N81").Value = -41746
$ws.Range("H82").Value = 1313.8723
$ws.Range("I82").Value = 882.04
$ws.Range("J82").Value = 1804.591
$ws.Range("K82").Value = 882.04
$ws.Range("L82").Value = 1804.591
$ws.Range("M82").Value = -521.04
$ws.Range("N82").Value = -2526.591
$ws.Range("H83").Value = 39750
$ws.Range("J83").Value = 39750
$ws.Range("L83").Value = 119250
$ws.Range("N83").Value = -130482
$ws.Range("H84").Value = 39750
$ws.Range("J84").Value = 39750
$ws.Range("L84").Value = 119250
$ws.Range("N84").Value = -129234
$ws.Range("H85").Value = 1313.8723
$ws.Range("I85").Value = 882.04
$ws.Range("J85").Value = 1804.591
$ws.Range("K85").Value = 882.04
$ws.Range("L85").Value = 1804.591
$ws.Range("M85").Value = 365.96
$ws.Range("N85").Value = -4300.591
$ws.Range("H86").Value = 34111.43
$ws.Range("J86").Value = 34111.43
$ws.Range("L86").Value = 34111.43
$ws.Range("N86").Value = -36483.43
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52246
$ws.Range("H89").Value = 34111.43
$ws.Range("J89").Value = 34111.43
$ws.Range("L89").Value = 102334.29
$ws.Range("N89").Value = -114190.29
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161232
$ws.Range("H127").Value = 40286.25
$ws.Range("J127").Value = 40286.25
$ws.Range("L127").Value = 40286.25
$ws.Range("N127").Value = -50206.25
$ws.Range("H136").Value = 1888.0741
$ws.Range("I136").Value = 1011.5
$ws.Range("K136").Value = 3034.5
$ws.Range("M136").Value = -484.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 28000
$ws.Range("J124").Value = 28000
$ws.Range("L124").Value = 28000
$ws.Range("N124").Value = -37820
$ws.Range("H132").Value = 6062157
$ws.Range("I132").Value = 1092.8096
$ws.Range("J132").Value = 25644058
$ws.Range("K132").Value = 3278.4288
$ws.Range("L132").Value = 76932174
$ws.Range("M132").Value = -748.4288000000001
$ws.Range("N132").Value = -76937234
